$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-13 16:48:27"
$ws.Range("E3").Value = "2026-02-13 16:48:30"
$ws.Range("I3").Value = "5.0 mm"
$ws.Range("K3").Value = "6.3 MJ/m2"
$ws.Range("E4").Value = "2026-02-13 16:48:32"
$ws.Range("H4").Value = "'73%"
$ws.Range("I4").Value = "4.0 mm"
$ws.Range("J4").Value = "996.2 hPa"
$ws.Range("K4").Value = "3.2 MJ/m2"
$ws.Range("E5").Value = "2026-02-13 16:48:35"
$ws.Range("K5").Value = "4.7 MJ/m2"
$ws.Range("E6").Value = "2026-02-13 16:48:37"
$ws.Range("H6").Value = "'72%"
$ws.Range("I6").Value = "2.4 mm"
$ws.Range("J6").Value = "996.3 hPa"
$ws.Range("K6").Value = "3.4 MJ/m2"
$ws.Range("E7").Value = "2026-02-13 16:48:40"
$ws.Range("H7").Value = "'69%"
$ws.Range("I7").Value = "14.7 mm"
$ws.Range("J7").Value = "996.4 hPa"
$ws.Range("K7").Value = "1.6 MJ/m2"
$ws.Range("N7").Value = "10.9 °C 16:26 TU"
$ws.Range("O7").Value = "13.2 °C"
$ws.Range("E8").Value = "2026-02-13 16:48:42"
$ws.Range("H8").Value = "'76%"
$ws.Range("I8").Value = "16.7 mm"
$ws.Range("J8").Value = "996.4 hPa"
$ws.Range("L8").Value = "45.4 km/h - 256º 16:29 TU"
$ws.Range("E9").Value = "2026-02-13 16:48:44"
$ws.Range("I9").Value = "0.9 mm"
$ws.Range("E10").Value = "2026-02-13 16:48:47"
$ws.Range("H10").Value = "'85%"
$ws.Range("I10").Value = "15.8 mm"
$ws.Range("K10").Value = "2.5 MJ/m2"
$ws.Range("E11").Value = "2026-02-13 16:48:49"
$ws.Range("I11").Value = "15.1 mm"
$ws.Range("E12").Value = "2026-02-13 16:48:51"
$ws.Range("H12").Value = "'82%"
$ws.Range("I12").Value = "3.4 mm"
$ws.Range("E13").Value = "2026-02-13 16:48:54"
$ws.Range("E14").Value = "2026-02-13 16:48:56"
$ws.Range("I14").Value = "17.2 mm"
$ws.Range("K14").Value = "1.5 MJ/m2"
$ws.Range("O14").Value = "10.6 °C"
$ws.Range("E15").Value = "2026-02-13 16:48:59"
$ws.Range("H15").Value = "'74%"
$ws.Range("I15").Value = "0.8 mm"
$ws.Range("E16").Value = "2026-02-13 16:49:01"
$ws.Range("H16").Value = "'79%"
$ws.Range("K16").Value = "6.1 MJ/m2"
$ws.Range("E17").Value = "2026-02-13 16:49:04"
$ws.Range("H17").Value = "'88%"
$ws.Range("I17").Value = "4.4 mm"
$ws.Range("K17").Value = "1.5 MJ/m2"
$ws.Range("E18").Value = "2026-02-13 16:49:06"
$ws.Range("H18").Value = "'80%"
$ws.Range("I18").Value = "6.3 mm"
$ws.Range("J18").Value = "996.4 hPa"
$ws.Range("E19").Value = "2026-02-13 16:49:09"
$ws.Range("I19").Value = "10.1 mm"
$ws.Range("O19").Value = "3.9 °C"
$ws.Range("E20").Value = "2026-02-13 16:49:11"
$ws.Range("I20").Value = "18.3 mm"
$ws.Range("K20").Value = "4.2 MJ/m2"
$ws.Range("O20").Value = "-3.9 °C"
$ws.Range("E21").Value = "2026-02-13 16:49:13"
$ws.Range("J21").Value = "999.1 hPa"
$ws.Range("E22").Value = "2026-02-13 16:49:16"
$ws.Range("H22").Value = "'91%"
$ws.Range("K22").Value = "6.7 MJ/m2"
$ws.Range("O22").Value = "-5.3 °C"
$ws.Range("E23").Value = "2026-02-13 16:49:18"
$ws.Range("I23").Value = "6.1 mm"
$ws.Range("K23").Value = "5.8 MJ/m2"
$ws.Range("E24").Value = "2026-02-13 16:49:21"
$ws.Range("I24").Value = "11.9 mm"
$ws.Range("J24").Value = "996.6 hPa"
$ws.Range("K24").Value = "1.5 MJ/m2"
$ws.Range("E25").Value = "2026-02-13 16:49:23"
$ws.Range("I25").Value = "8.3 mm"
$ws.Range("E26").Value = "2026-02-13 16:49:26"
$ws.Range("E27").Value = "2026-02-13 16:49:28"
$ws.Range("H27").Value = "'81%"
$ws.Range("I27").Value = "17.6 mm"
$ws.Range("K27").Value = "3.4 MJ/m2"
$ws.Range("E28").Value = "2026-02-13 16:49:31"
$ws.Range("I28").Value = "5.4 mm"
$ws.Range("J28").Value = "996.7 hPa"
$ws.Range("K28").Value = "3.0 MJ/m2"
$ws.Range("E29").Value = "2026-02-13 16:49:33"
$ws.Range("I29").Value = "11.0 mm"
$ws.Range("E30").Value = "2026-02-13 16:49:36"
$ws.Range("H30").Value = "'75%"
$ws.Range("I30").Value = "0.9 mm"
$ws.Range("J30").Value = "996.1 hPa"
$ws.Range("O30").Value = "9.2 °C"
$ws.Range("E31").Value = "2026-02-13 16:49:38"
$ws.Range("H31").Value = "'69%"
$ws.Range("I31").Value = "0.5 mm"
$ws.Range("J31").Value = "995.1 hPa"
$ws.Range("N31").Value = "8.8 °C 16:02 TU"
$ws.Range("O31").Value = "10.6 °C"
$ws.Range("E32").Value = "2026-02-13 16:49:41"
$ws.Range("I32").Value = "22.1 mm"
$ws.Range("K32").Value = "1.2 MJ/m2"
$ws.Range("E33").Value = "2026-02-13 16:49:43"
$ws.Range("J33").Value = "998.1 hPa"
$ws.Range("E34").Value = "2026-02-13 16:49:46"
$ws.Range("G34").Value = "102 cm"
$ws.Range("I34").Value = "7.9 mm"
$ws.Range("K34").Value = "4.5 MJ/m2"
$ws.Range("E35").Value = "2026-02-13 16:49:48"
$ws.Range("H35").Value = "'73%"
$ws.Range("J35").Value = "996.6 hPa"
$ws.Range("K35").Value = "2.7 MJ/m2"
$ws.Range("E36").Value = "2026-02-13 16:49:51"
$ws.Range("I36").Value = "6.7 mm"
$ws.Range("J36").Value = "996.4 hPa"
$ws.Range("O36").Value = "10.8 °C"
$ws.Range("E37").Value = "2026-02-13 16:49:53"
$ws.Range("I37").Value = "11.3 mm"
$ws.Range("J37").Value = "998.2 hPa"
$ws.Range("O37").Value = "3.6 °C"
$ws.Range("E38").Value = "2026-02-13 16:49:56"
$ws.Range("I38").Value = "10.3 mm"
$ws.Range("E39").Value = "2026-02-13 16:49:58"
$ws.Range("I39").Value = "14.2 mm"
$ws.Range("K39").Value = "6.9 MJ/m2"
$ws.Range("E40").Value = "2026-02-13 16:50:01"
$ws.Range("J40").Value = "999.6 hPa"
$ws.Range("E41").Value = "2026-02-13 16:50:03"
$ws.Range("H41").Value = "'77%"
$ws.Range("J41").Value = "996.3 hPa"
$ws.Range("K41").Value = "1.5 MJ/m2"
$ws.Range("E42").Value = "2026-02-13 16:50:06"
$ws.Range("I42").Value = "7.0 mm"
$ws.Range("E43").Value = "2026-02-13 16:50:08"
$ws.Range("H43").Value = "'86%"
$ws.Range("I43").Value = "11.0 mm"
$ws.Range("K43").Value = "2.2 MJ/m2"
$ws.Range("E44").Value = "2026-02-13 16:50:10"
$ws.Range("I44").Value = "1.8 mm"
$ws.Range("K44").Value = "5.5 MJ/m2"
$ws.Range("E45").Value = "2026-02-13 16:50:13"
$ws.Range("J45").Value = "994.3 hPa"
$ws.Range("E46").Value = "2026-02-13 16:50:15"
$ws.Range("J46").Value = "996.7 hPa"
$ws.Range("K46").Value = "1.7 MJ/m2"
$ws.Range("O46").Value = "8.6 °C"
